# Updated cryptos list on Tue Jun  4 15:53:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so values like "666.34" or
# "1.00" are not auto-converted to numbers (the source data stores every
# price/volume cell as a plain string, even when it looks numeric).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.204.95"
$ws.Range("E2").Value = "  +1.41%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.802.82"
$ws.Range("E3").Value = "  +0.42%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "666.34"
$ws.Range("E5").Value = "  +6.21%  "

# Row 6 - Solana
$ws.Range("D6").Value = "167.01"
$ws.Range("E6").Value = "  +1.18%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.801.12"

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.07%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  +1.51%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.66%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.76%  "

# Row 12 - Toncoin
$ws.Range("D12").Value = "6.99"
$ws.Range("E12").Value = "  +4.83%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -2.49%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "35.73"
$ws.Range("E14").Value = "  +0.34%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.441.49"
$ws.Range("E15").Value = "  +0.42%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.796.89"
$ws.Range("E16").Value = "  +0.67%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "70.117.84"
$ws.Range("E17").Value = "  +1.28%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "17.77"
$ws.Range("E18").Value = "  -0.92%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "7.18"
$ws.Range("E19").Value = "  +0.65%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.79%  "

# Row 21 - now Uniswap (was BitcoinCash)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "10.32"
$ws.Range("E21").Value = "  +6.94%  "

# Row 22 - now BitcoinCash (was Uniswap)
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "474.36"
$ws.Range("E22").Value = "  +1.16%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.714"
$ws.Range("E23").Value = "  +1.12%  "

# Row 24 - PEPE
$ws.Range("D24").Value = "0.0000146"
$ws.Range("E24").Value = "  -3.75%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "82.64"
$ws.Range("E25").Value = "  -0.65%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("E26").Value = "  +1.28%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "10.39"
$ws.Range("E27").Value = "  +3.66%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  -1.81%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - WrappedeETH
$ws.Range("D30").Value = "3.951.46"
$ws.Range("E30").Value = "  +0.43%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "2.81"
$ws.Range("E31").Value = "  +4.90%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +3.29%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "7.36"
$ws.Range("E33").Value = "  +0.50%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "29.18"
$ws.Range("E34").Value = "  +0.77%  "

# Row 35 - Kaspa
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  +17.44%  "

# Row 36 - Binance-PegBSC-USD
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "3.756.15"
$ws.Range("E37").Value = "  +0.54%  "

# Row 38 - Aptos
$ws.Range("E38").Value = "  +0.07%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  -0.62%  "

# Row 40 - now dogwifhat (was Filecoin)
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "3.33"
$ws.Range("E40").Value = "  -0.89%  "

# Row 41 - now Filecoin (was dogwifhat)
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "5.92"
$ws.Range("E41").Value = "  +1.61%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.965"
$ws.Range("E42").Value = "  -0.33%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  -0.03%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  -0.01%  "

# Row 45 - Stacks
$ws.Range("D45").Value = "2.07"
$ws.Range("E45").Value = "  +6.38%  "

# Row 46 - Arweave
$ws.Range("D46").Value = "45.01"
$ws.Range("E46").Value = "  +4.56%  "

# Row 47 - Monero
$ws.Range("D47").Value = "158.73"
$ws.Range("E47").Value = "  +3.72%  "

# Row 48 - OKB
$ws.Range("D48").Value = "47.74"
$ws.Range("E48").Value = "  +1.97%  "

# Row 49 - TheGraph
$ws.Range("E49").Value = "  +0.11%  "

# Row 50 - ONDO
$ws.Range("E50").Value = "  +1.94%  "

# Row 51 - Cosmos
$ws.Range("D51").Value = "8.50"
$ws.Range("E51").Value = "  +0.73%  "

# Restore the default (unstyled) cell style for the Price column now that
# the text values are locked in, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"

Write-Output "cryptos updated"
